$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.274.65'
$ws.Range("E2").Value = '  +5.26%  '
$ws.Range("D3").Value = '''2.296.49'
$ws.Range("E3").Value = '  +5.59%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''252.58'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Value = '''0.643'
$ws.Range("E6").Value = '  +5.35%  '
$ws.Range("D7").Value = '''73.31'
$ws.Range("E7").Value = '  +10.87%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '''0.668'
$ws.Range("E9").Value = '  +13.84%  '
$ws.Range("D10").Value = '''39.80'
$ws.Range("E10").Value = '  +9.76%  '
$ws.Range("D11").Value = '''0.0984'
$ws.Range("E11").Value = '  +5.49%  '
$ws.Range("D12").Value = '''59.98'
$ws.Range("E12").Value = '  +2.08%  '
$ws.Range("D13").Value = '''7.53'
$ws.Range("E13").Value = '  +9.97%  '
$ws.Range("D14").Value = '''0.105'
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("D15").Value = '''2.628.88'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("D16").Value = '''15.20'
$ws.Range("E16").Value = '  +6.83%  '
$ws.Range("D17").Value = '''0.899'
$ws.Range("E17").Value = '  +6.74%  '
$ws.Range("D18").Value = '''2.290.08'
$ws.Range("E18").Value = '  +5.10%  '
$ws.Range("D19").Value = '''43.177.40'
$ws.Range("E19").Value = '  +5.18%  '
$ws.Range("E20").Value = '  +7.44%  '
$ws.Range("D21").Value = '''6.42'
$ws.Range("E21").Value = '  +6.56%  '
$ws.Range("D22").Value = '''73.82'
$ws.Range("E22").Value = '  +3.31%  '
$ws.Range("D23").Value = '''238.47'
$ws.Range("E23").Value = '  +3.86%  '
$ws.Range("D24").Value = '''2.20'
$ws.Range("E24").Value = '  +8.60%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '''11.93'
$ws.Range("E25").Value = '  +6.36%  '
$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").Value = '''3.91'
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = '''2.46'
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("D30").Value = '''2.20'
$ws.Range("E30").Value = '  +8.49%  '
$ws.Range("D31").Value = '''168.37'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").Value = '''21.29'
$ws.Range("E32").Value = '  +5.78%  '
$ws.Range("D33").Value = '''6.33'
$ws.Range("E33").Value = '  +12.53%  '
$ws.Range("D34").Value = '''0.129'
$ws.Range("E34").Value = '  +7.55%  '
$ws.Range("D35").Value = '''0.0813'
$ws.Range("E35").Value = '  +9.38%  '
$ws.Range("D36").Value = '''31.34'
$ws.Range("E36").Value = '  +27.69%  '
$ws.Range("D37").Value = '''4.84'
$ws.Range("E37").Value = '  +22.34%  '
$ws.Range("E38").Value = '  +5.36%  '
$ws.Range("D39").Value = '''4.82'
$ws.Range("E39").Value = '  +6.94%  '
$ws.Range("D40").Value = '''0.0313'
$ws.Range("E40").Value = '  +3.08%  '
$ws.Range("D41").Value = '''13.63'
$ws.Range("E41").Value = '  +21.17%  '
$ws.Range("E42").Value = '  +6.72%  '
$ws.Range("D43").Value = '''6.15'
$ws.Range("E43").Value = '  +11.77%  '
$ws.Range("D44").Value = '''0.214'
$ws.Range("E44").Value = '  +14.14%  '
$ws.Range("D45").Value = '''9.24'
$ws.Range("E45").Value = '  +9.06%  '
$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").Value = '''62.29'
$ws.Range("E46").Value = '  +2.92%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '''4.94'
$ws.Range("E47").Value = '  -9.42%  '
$ws.Range("E48").Value = '  +5.49%  '
$ws.Range("D49").Value = '''1.19'
$ws.Range("E49").Value = '  +5.32%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.20'
$ws.Range("E51").Value = '  +5.86%  '
